$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 104 (shifts existing rows 104:201 down to 105:202,
# and the sheet dimension grows from A1:R201 to A1:R202).
$ws.Rows.Item(104).Insert()

# Populate the newly inserted row with the new "Acelga / Feria Lagunitas de
# Puerto Montt / Los Lagos" price record.
$ws.Cells.Item(104, 1).Value  = 4
$ws.Cells.Item(104, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(104, 3).Value  = "Los Lagos"
$ws.Cells.Item(104, 4).Value  = 44778
$ws.Cells.Item(104, 5).Value  = 10
$ws.Cells.Item(104, 6).Value  = 100112009
$ws.Cells.Item(104, 7).Value  = "Acelga"
$ws.Cells.Item(104, 8).Value  = "Sin especificar"
$ws.Cells.Item(104, 9).Value  = "Primera"
$ws.Cells.Item(104, 10).Value = 200
$ws.Cells.Item(104, 11).Value = 1200
$ws.Cells.Item(104, 12).Value = 1500
$ws.Cells.Item(104, 13).Value = 1350
$ws.Cells.Item(104, 14).Value = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(104, 15).Value = "Región de Los Lagos"
$ws.Cells.Item(104, 16).Value = 900
$ws.Cells.Item(104, 17).Value = 1.5
$ws.Cells.Item(104, 18).Value = "Hortaliza"
